$d = $word.ActiveDocument

# Insert a new paragraph right after the last paragraph in the document
# (the final bullet "The location wherein an opportunity exist").
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last

# The new paragraph must NOT be part of the bulleted list - strip the
# list numbering / ListParagraph style that would otherwise be inherited
# from the previous paragraph, and reset it back to the Normal style.
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = "Normal"

# Match the size formatting (28 half-points => 14pt) used throughout the
# rest of the document, then set the paragraph text.
$newPara.Range.Font.Size = 14
$newPara.Range.Text = "Target Audience is XYZ leadership."
